$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usages")

# Insert a new row at 49 (pushes "Limit Switches" and everything below
# down by one). Excel's row Insert shifts cells down and carries the
# formatting of the row above along with it, matching the "IR Array"
# group (rows 42-48) that this new row belongs to.
$ws.Rows.Item(49).Insert()

# Add the missing "IR Array" pin: D43.
$ws.Cells.Item(49, 1).Value = "IR Array"
$ws.Cells.Item(49, 3).Value = "D43"

$wb.Save()
